$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.666.71'
$ws.Range('D3').Value = '1.953.27'
$ws.Range('E3').Value = '  +1.90%  '
$ws.Range('E4').Value = '  -0.24%  '
$ws.Range('D5').Value = "'244.48"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.27%  '
$ws.Range('D6').Value = "'0.611"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.04%  '
$ws.Range('D7').Value = "'58.43"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +6.20%  '
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('D9').Value = "'0.368"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.15%  '
$ws.Range('D10').Value = "'0.0813"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.79%  '
$ws.Range('E11').Value = '  +0.47%  '
$ws.Range('D12').Value = "'22.38"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +10.00%  '
$ws.Range('D13').Value = '2.237.70'
$ws.Range('E13').Value = '  +1.69%  '
$ws.Range('D14').Value = "'0.819"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.84%  '
$ws.Range('D15').Value = "'13.58"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.58%  '
$ws.Range('D16').Value = "'5.24"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.70%  '
$ws.Range('D17').Value = '1.953.49'
$ws.Range('E17').Value = '  +3.12%  '
$ws.Range('D18').Value = '36.625.69'
$ws.Range('E18').Value = '  +1.96%  '
$ws.Range('D19').Value = "'69.71"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.50%  '
$ws.Range('D20').Value = '0.0₃0861'
$ws.Range('E20').Value = '  +1.56%  '
$ws.Range('D21').Value = "'228.59"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.89%  '
$ws.Range('D22').Value = "'5.04"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.59%  '
$ws.Range('E23').Value = '  -0.22%  '
$ws.Range('E24').Value = '  +1.18%  '
$ws.Range('E25').Value = '  +4.06%  '
$ws.Range('D26').Value = "'9.27"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.54%  '
$ws.Range('D27').Value = "'160.60"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.87%  '
$ws.Range('D28').Value = "'0.134"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +17.12%  '
$ws.Range('D29').Value = "'19.34"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.25%  '
$ws.Range('E30').Value = '  +2.77%  '
$ws.Range('E31').Value = '  +0.75%  '
$ws.Range('D32').Value = "'4.69"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.11%  '
$ws.Range('D33').Value = "'0.0621"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.64%  '
$ws.Range('E34').Value = '  +0.96%  '
$ws.Range('D35').Value = "'6.25"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +6.79%  '
$ws.Range('E36').Value = '  -0.22%  '
$ws.Range('D37').Value = "'3.42"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +23.18%  '
$ws.Range('E38').Value = '  +5.52%  '
$ws.Range('E39').Value = '  -1.24%  '
$ws.Range('D40').Value = "'0.100"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.54%  '
$ws.Range('E41').Value = '  +1.84%  '
$ws.Range('E42').Value = '  +3.40%  '
$ws.Range('E43').Value = '  +1.06%  '
$ws.Range('D44').Value = "'16.05"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +5.61%  '
$ws.Range('E45').Value = '  +3.45%  '
$ws.Range('D46').Value = '1.348.76'
$ws.Range('E46').Value = '  +2.07%  '
$ws.Range('D47').Value = "'87.36"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.84%  '
$ws.Range('D48').Value = "'7.20"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.14%  '
$ws.Range('E49').Value = '  +1.75%  '
$ws.Range('D50').Value = '2.130.04'
$ws.Range('E50').Value = '  +1.72%  '
$ws.Range('D51').Value = "'43.54"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.54%  '
